$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.677.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.49%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.536.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3925"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.80%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3160"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.35%  "

# Row 12
$ws.Range("E12").Value = "  +0.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.610"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.49%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.600"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.539.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001095"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06583"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.42%  "

# Row 20
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.122"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.392"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.677.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.341"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.39%  "

# Row 29
$ws.Range("E29").Value = "  -0.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.711.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.50%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.880"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.62%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9592"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.97%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08150"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.633"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06060"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.097"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.70%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2021"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.40%  "

# Row 40
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.434"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.13%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.181"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.90%  "

# Row 42
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
$ws.Range("E43").Value = "  -1.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5714"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.732"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.20%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5466"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.163"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.77%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.860"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06688"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.93%  "
